$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 370.6875
$ws.Range("J2").Value = 674.2
$ws.Range("L2").Value = 674.2
$ws.Range("N2").Value = -900.2

$ws.Range("H51").Value = 3391.6667
$ws.Range("J51").Value = 3662.5
$ws.Range("L51").Value = 3662.5
$ws.Range("N51").Value = -4630.5

$ws.Range("H107").Value = 346.57895
$ws.Range("I107").Value = 117
$ws.Range("J107").Value = 989.4
$ws.Range("K107").Value = 117
$ws.Range("L107").Value = 989.4
$ws.Range("M107").Value = 1803
$ws.Range("N107").Value = -4829.4

$ws.Range("H116").Value = 4585.7144
$ws.Range("I116").Value = 4700.2
$ws.Range("K116").Value = 4700.2
$ws.Range("M116").Value = -1258.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 144.55556
$ws.Range("I4").Value = 176.28572
$ws.Range("J4").Value = 33.5
$ws.Range("K4").Value = 176.28572
$ws.Range("L4").Value = 33.5
$ws.Range("M4").Value = -60.28572
$ws.Range("N4").Value = -265.5

$ws.Range("H5").Value = 102.166664
$ws.Range("I5").Value = 102.166664
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 102.166664
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 9.833336000000003
$ws.Range("N5").ClearContents()

$ws.Range("H6").Value = 17500000
$ws.Range("I6").Value = 17500000
$ws.Range("K6").Value = 17500000
$ws.Range("M6").Value = -17499827

$ws.Range("H32").Value = 5363.647
$ws.Range("I32").Value = 4995.879
$ws.Range("K32").Value = 4995.879
$ws.Range("M32").Value = -4708.879

$ws.Range("H41").Value = 1331.25
$ws.Range("I41").Value = 1331.25
$ws.Range("K41").Value = 1331.25
$ws.Range("M41").Value = -917.25

$ws.Range("H44").Value = 12130.5
$ws.Range("J44").Value = 12130.5
$ws.Range("L44").Value = 12130.5
$ws.Range("N44").Value = -13106.5

$ws.Range("H45").Value = 3347.1177
$ws.Range("I45").Value = 2778.6428
$ws.Range("K45").Value = 2778.6428
$ws.Range("M45").Value = -2401.6428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 102.166664
$ws.Range("I4").Value = 102.166664
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 102.166664
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 12.833336
$ws.Range("N4").ClearContents()

$ws.Range("H22").Value = 401.8
$ws.Range("I22").Value = 292.69232
$ws.Range("J22").Value = 1111
$ws.Range("K22").Value = 292.69232
$ws.Range("L22").Value = 1111
$ws.Range("M22").Value = -119.69232
$ws.Range("N22").Value = -1457

$ws.Range("H80").Value = 222.44444
$ws.Range("I80").Value = 179.6
$ws.Range("K80").Value = 179.6
$ws.Range("M80").Value = 818.4

$ws.Range("H83").Value = 222.44444
$ws.Range("I83").Value = 179.6
$ws.Range("K83").Value = 898
$ws.Range("M83").Value = 4094

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1549.5
$ws.Range("I22").Value = 216
$ws.Range("J22").Value = 2549.625
$ws.Range("K22").Value = 216
$ws.Range("L22").Value = 2549.625
$ws.Range("M22").Value = 134
$ws.Range("N22").Value = -3249.625

$ws.Range("H31").Value = 4145.4585
$ws.Range("J31").Value = 9920.571
$ws.Range("L31").Value = 9920.571
$ws.Range("N31").Value = -10510.571

$ws.Range("H34").Value = 4145.4585
$ws.Range("J34").Value = 9920.571
$ws.Range("L34").Value = 9920.571
$ws.Range("N34").Value = -10324.571

$ws.Range("H105").Value = 3651.5
$ws.Range("I105").Value = 3081
$ws.Range("K105").Value = 3081
$ws.Range("M105").Value = -1334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.23077
$ws.Range("I2").Value = 33.615383
$ws.Range("J2").Value = 32.846153
$ws.Range("K2").Value = 201.692298
$ws.Range("L2").Value = 197.076918
$ws.Range("M2").Value = -88.69229799999999
$ws.Range("N2").Value = -423.076918

$ws.Range("H14").Value = 4181.2
$ws.Range("I14").Value = 4181.2
$ws.Range("K14").Value = 12543.6
$ws.Range("M14").Value = -12370.6

$ws.Range("H34").Value = 4020
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15168

$ws.Range("H39").Value = 7108.4614
$ws.Range("I39").Value = 638
$ws.Range("J39").Value = 9049.6
$ws.Range("K39").Value = 1914
$ws.Range("L39").Value = 27148.8
$ws.Range("M39").Value = -1620
$ws.Range("N39").Value = -27736.8

$ws.Range("H55").Value = 4475
$ws.Range("I55").Value = 1016.6667
$ws.Range("J55").Value = 6550
$ws.Range("K55").Value = 3050.0001
$ws.Range("L55").Value = 19650
$ws.Range("M55").Value = -2873.0001
$ws.Range("N55").Value = -20004

$ws.Range("H57").Value = 1726.5
$ws.Range("I57").Value = 975
$ws.Range("J57").Value = 2478
$ws.Range("K57").Value = 2925
$ws.Range("L57").Value = 7434
$ws.Range("M57").Value = -2366
$ws.Range("N57").Value = -8552

$ws.Range("H92").Value = 6237.25
$ws.Range("J92").Value = 5149.6665
$ws.Range("L92").Value = 15448.9995
$ws.Range("N92").Value = -17944.9995

$ws.Range("H98").Value = 416
$ws.Range("I98").Value = 229
$ws.Range("J98").Value = 790
$ws.Range("K98").Value = 687
$ws.Range("L98").Value = 2370
$ws.Range("M98").Value = 811
$ws.Range("N98").Value = -5366

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5318.4287
$ws.Range("I100").Value = 1898.375
$ws.Range("K100").Value = 1898.375
$ws.Range("M100").Value = -1357.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H74").Value = 16480.5
$ws.Range("I74").Value = 4992
$ws.Range("K74").Value = 4992
$ws.Range("M74").Value = -4056

$ws.Range("H77").Value = 16480.5
$ws.Range("I77").Value = 4992
$ws.Range("K77").Value = 14976
$ws.Range("M77").Value = -10296

$ws.Range("H81").Value = 999
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 999
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H96").Value = 993.25
$ws.Range("I96").Value = 900
$ws.Range("J96").Value = 1024.3334
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 1024.3334
$ws.Range("M96").Value = 473
$ws.Range("N96").Value = -3770.3334

$ws.Range("H100").Value = 477.7857
$ws.Range("I100").Value = 286.125
$ws.Range("K100").Value = 572.25
$ws.Range("M100").Value = -31.25
